$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data row (row 2) with the new values.
# A leading apostrophe keeps the cells' existing "stored as text"
# (quote-prefix) formatting while the apostrophe itself is not
# persisted as part of the stored string.
$ws.Range("A2").Value = "btorres"
$ws.Range("F2").Value = "'ACHACALTANAS1"
$ws.Range("D2").Value = "'AA21180FH5N8"
$ws.Range("C2").Value = "'PRESTAMO YA"

# Clear the former Estado / Transaccion / Fecha values for row 2
$ws.Range("G2:I2").ClearContents()

# Update the active selection as recorded in the saved workbook
$ws.Range("J3").Select() | Out-Null
